$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows at position 932 (pushes existing rows 932:1019 down to 934:1021)
$ws.Rows("932:933").Insert()

# --- New row 932 ---
$ws.Range("A932").Value = 5
$ws.Range("B932").Value = "Macroferia Regional de Talca"
$ws.Range("C932").Value = "Maule"
$ws.Range("D932").Value = 45166
$ws.Range("E932").Value = 7
$ws.Range("F932").Value = 100112004
$ws.Range("G932").Value = "Cebolla"
$ws.Range("H932").Value = "Sin especificar"
$ws.Range("I932").Value = "1a (guarda)"
$ws.Range("J932").Value = 2000
$ws.Range("K932").Value = 19000
$ws.Range("L932").Value = 19000
$ws.Range("M932").Value = 19000
$ws.Range("N932").Value = "`$/malla 25 kilos"
$ws.Range("O932").Value = "Región del Maule"
$ws.Range("P932").Value = 760
$ws.Range("Q932").Value = 25
$ws.Range("R932").Value = "Hortaliza"

# --- New row 933 ---
$ws.Range("A933").Value = 5
$ws.Range("B933").Value = "Macroferia Regional de Talca"
$ws.Range("C933").Value = "Maule"
$ws.Range("D933").Value = 45166
$ws.Range("E933").Value = 7
$ws.Range("F933").Value = 100112004
$ws.Range("G933").Value = "Cebolla"
$ws.Range("H933").Value = "Sin especificar"
$ws.Range("I933").Value = "2a (guarda)"
$ws.Range("J933").Value = 800
$ws.Range("K933").Value = 18000
$ws.Range("L933").Value = 18000
$ws.Range("M933").Value = 18000
$ws.Range("N933").Value = "`$/malla 25 kilos"
$ws.Range("O933").Value = "Región del Maule"
$ws.Range("P933").Value = 720
$ws.Range("Q933").Value = 25
$ws.Range("R933").Value = "Hortaliza"
